# Natmi following Dr Hou advice
# Extend the LR-pairs table from 2 target/sending clusters (FAPs, sCs) to 3
# (ECs, FAPs, sCs), which grows the data block from rows 2-5 to rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Ntf3", "Ntrk3", "FAPs", 2, 0.6666666666666666, 1.788130666666667, 5.364392, 0.08755714261138148, 0.08755714261138148, 3, 1, 1.467507333333333, 4.402521999999999, 0.1890754490804, 0.1890754490804, 2.624094866291555, 23.616853796624, 0.01655490605944358, 0.01655490605944358),
    @("ECs", "Ntf3", "Ntrk3", "sCs", 2, 0.6666666666666666, 1.788130666666667, 5.364392, 0.08755714261138148, 0.08755714261138148, 3, 1, 6.293983333333333, 18.88195, 0.8109245509196, 0.8109245509195999, 11.25446461382222, 101.2901815244, 0.07100223655193789, 0.07100223655193789),
    @("FAPs", "Ntf3", "Ntrk3", "FAPs", 3, 1, 11.451921, 34.355763, 0.560751794520949, 0.560751794520949, 3, 1, 1.467507333333333, 4.402521999999999, 0.1890754490804, 0.1890754490804, 16.805778048254, 151.252002434286, 0.1060243973716886, 0.1060243973716886),
    @("FAPs", "Ntf3", "Ntrk3", "sCs", 3, 1, 11.451921, 34.355763, 0.560751794520949, 0.560751794520949, 3, 1, 6.293983333333333, 18.88195, 0.8109245509196, 0.8109245509195999, 72.07819990865001, 648.7037991778501, 0.4547273971492604, 0.4547273971492604),
    @("sCs", "Ntf3", "Ntrk3", "FAPs", 3, 1, 7.182390333333333, 21.547171, 0.3516910628676694, 0.3516910628676694, 3, 1, 1.467507333333333, 4.402521999999999, 0.1890754490804, 0.1890754490804, 10.54021048502911, 94.86189436526197, 0.06649614564926778, 0.06649614564926778),
    @("sCs", "Ntf3", "Ntrk3", "sCs", 3, 1, 7.182390333333333, 21.547171, 0.3516910628676694, 0.3516910628676694, 3, 1, 6.293983333333333, 18.88195, 0.8109245509196, 0.8109245509195999, 45.20584505149444, 406.8526054634499, 0.2851949172184017, 0.2851949172184016)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
